# Auto-generated Excel COM-interop script applying the meteocat daily-summary refresh.
# Commit: "Update automatic: dades i banners [2026-02-09 06:20]"
# Refreshes DATA_EXTRACCIO scrape timestamps and the weather-derived metrics that shifted
# between the 05:48-05:50 run and the 06:18-06:20 run, for every station row (2-46).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values that are pure percentages (e.g. "78%") must be forced to stay as literal
# text -- like the source workbook stores them -- instead of being auto-converted by
# Excel into a numeric percentage value when assigned through .Value.
$percentCells = @(
    'H8'
    'H9'
    'H11'
    'H12'
    'H15'
    'H16'
    'H17'
    'H23'
    'H25'
    'H30'
    'H33'
    'H34'
    'H36'
    'H39'
)

$updates = [ordered]@{
    'E2' = '2026-02-09 06:18:41'
    'E3' = '2026-02-09 06:18:44'
    'O3' = '-6.1 °C'
    'E4' = '2026-02-09 06:18:47'
    'O4' = '4.2 °C'
    'E5' = '2026-02-09 06:18:50'
    'M5' = '-4.1 °C 5:45 TU'
    'E6' = '2026-02-09 06:18:52'
    'O6' = '6.3 °C'
    'E7' = '2026-02-09 06:18:55'
    'J7' = '1008.3 hPa'
    'E8' = '2026-02-09 06:18:58'
    'H8' = '78%'
    'E9' = '2026-02-09 06:19:01'
    'H9' = '89%'
    'O9' = '6.4 °C'
    'E10' = '2026-02-09 06:19:03'
    'N10' = '2.4 °C 5:37 TU'
    'O10' = '5.2 °C'
    'E11' = '2026-02-09 06:19:06'
    'H11' = '97%'
    'E12' = '2026-02-09 06:19:09'
    'H12' = '92%'
    'N12' = '3.8 °C 5:54 TU'
    'O12' = '7.0 °C'
    'E13' = '2026-02-09 06:19:11'
    'O13' = '-1.9 °C'
    'E14' = '2026-02-09 06:19:14'
    'N14' = '6.2 °C 5:54 TU'
    'O14' = '7.2 °C'
    'E15' = '2026-02-09 06:19:16'
    'H15' = '88%'
    'N15' = '2.5 °C 5:36 TU'
    'O15' = '5.2 °C'
    'E16' = '2026-02-09 06:19:19'
    'H16' = '63%'
    'I16' = '0.1 mm'
    'N16' = '-6.5 °C 5:51 TU'
    'O16' = '-5.3 °C'
    'E17' = '2026-02-09 06:19:22'
    'H17' = '96%'
    'N17' = '-0.8 °C 5:59 TU'
    'E18' = '2026-02-09 06:19:24'
    'N18' = '3.5 °C 5:30 TU'
    'O18' = '5.9 °C'
    'E19' = '2026-02-09 06:19:27'
    'L19' = '13.7 km/h - 226º 5:59 TU'
    'N19' = '2.7 °C 5:56 TU'
    'E20' = '2026-02-09 06:19:30'
    'M20' = '-5.0 °C 5:58 TU'
    'O20' = '-6.3 °C'
    'E21' = '2026-02-09 06:19:32'
    'J21' = '1009.9 hPa'
    'E22' = '2026-02-09 06:19:35'
    'N22' = '-8.1 °C 5:52 TU'
    'E23' = '2026-02-09 06:19:38'
    'H23' = '84%'
    'O23' = '-5.8 °C'
    'E24' = '2026-02-09 06:19:41'
    'J24' = '1009.3 hPa'
    'E25' = '2026-02-09 06:19:43'
    'H25' = '75%'
    'O25' = '-4.3 °C'
    'E26' = '2026-02-09 06:19:46'
    'E27' = '2026-02-09 06:19:48'
    'E28' = '2026-02-09 06:19:51'
    'O28' = '3.4 °C'
    'E29' = '2026-02-09 06:19:54'
    'N29' = '2.6 °C 5:34 TU'
    'O29' = '5.1 °C'
    'E30' = '2026-02-09 06:19:56'
    'H30' = '95%'
    'N30' = '4.8 °C 5:59 TU'
    'O30' = '6.3 °C'
    'E31' = '2026-02-09 06:19:58'
    'O31' = '8.9 °C'
    'E32' = '2026-02-09 06:20:01'
    'E33' = '2026-02-09 06:20:04'
    'H33' = '95%'
    'J33' = '1009.9 hPa'
    'O33' = '-0.8 °C'
    'E34' = '2026-02-09 06:20:07'
    'H34' = '73%'
    'O34' = '-3.4 °C'
    'E35' = '2026-02-09 06:20:09'
    'N35' = '3.2 °C 5:30 TU'
    'E36' = '2026-02-09 06:20:12'
    'H36' = '85%'
    'N36' = '4.6 °C 5:46 TU'
    'O36' = '8.0 °C'
    'E37' = '2026-02-09 06:20:15'
    'J37' = '1009.5 hPa'
    'N37' = '1.2 °C 5:30 TU'
    'O37' = '3.0 °C'
    'E38' = '2026-02-09 06:20:17'
    'K38' = '-0.1 MJ/m2'
    'N38' = '3.7 °C 5:59 TU'
    'O38' = '5.8 °C'
    'E39' = '2026-02-09 06:20:20'
    'H39' = '82%'
    'K39' = '-0.1 MJ/m2'
    'E40' = '2026-02-09 06:20:23'
    'E41' = '2026-02-09 06:20:25'
    'O41' = '11.0 °C'
    'E42' = '2026-02-09 06:20:28'
    'N42' = '3.7 °C 5:47 TU'
    'O42' = '6.1 °C'
    'E43' = '2026-02-09 06:20:31'
    'L43' = '17.6 km/h - 243º 5:34 TU'
    'N43' = '5.6 °C 5:59 TU'
    'O43' = '6.3 °C'
    'E44' = '2026-02-09 06:20:34'
    'O44' = '-7.1 °C'
    'E45' = '2026-02-09 06:20:36'
    'J45' = '1009.6 hPa'
    'M45' = '1.1 °C 5:50 TU'
    'O45' = '0.1 °C'
    'E46' = '2026-02-09 06:20:39'
}

foreach ($cellRef in $updates.Keys) {
    $range = $ws.Range($cellRef)
    if ($percentCells -contains $cellRef) {
        $range.NumberFormat = "@"
        $range.Value = $updates[$cellRef]
        $range.NumberFormat = "General"
    } else {
        $range.Value = $updates[$cellRef]
    }
}

